$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 5)
$ws.Range("A5").Value = "Padi & Sons"
$ws.Range("B5").Value = "Michael Jordan"
$ws.Range("C5").Value = "Acquire Los Angeles Lakers"
$ws.Range("D5").Value = 5555
$ws.Range("E5").Value = 2222
$ws.Range("G5").Value = 2222
$ws.Range("H5").Value = "United States"
$ws.Range("I5").Value = 43953
# Reuse the same date format already used by the Deposit Date column (I2:I4)
$ws.Range("I5").NumberFormat = "[$-F800]dddd\,\ mmmm\ dd\,\ yyyy"

# Column C narrows to fit the new (shorter) longest text in that column
$ws.Columns.Item(3).ColumnWidth = 52.16666666666667

# Window/view state: zoom in and move the selection to below the new row
$win = $ws.Application.ActiveWindow
$win.Zoom = 268
$ws.Range("A6").Select() | Out-Null
